$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

# Insert a new row at 36, shifting existing row 36 (and below) down by one.
$ws.Rows("36:36").Insert([Microsoft.Office.Interop.Excel.XlInsertShiftDirection]::xlShiftDown)

# Populate the newly inserted row with the new skill entry.
$ws.Cells.Item(36, 2).Value = "OAuth2/OpenID Connect"
$ws.Cells.Item(36, 3).Value = 3

# Match formatting used by the other skill-level cells in column C (style index 3 = centered).
$ws.Cells.Item(36, 3).HorizontalAlignment = [Microsoft.Office.Interop.Excel.XlHAlign]::xlHAlignCenter

# Update the saved view state to match where the author left the cursor/scroll.
$ws.Range("G35").Select()
$excel.ActiveWindow.ScrollRow = 25
$excel.ActiveWindow.ScrollColumn = 1
